$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

$ws.Range("A284").Value = "cs"
$ws.Range("B284").Value = "lab.build.button.create"
$ws.Range("C284").Value = "Nový build"

$ws.Range("A285").Value = "cs"
$ws.Range("B285").Value = "lab.build.button.list"
$ws.Range("C285").Value = "Seznam buildů"

$ws.Range("A286").Value = "cs"
$ws.Range("C286").Value = "Nejnovější buildy"
$ws.Range("B286").Value = "lab.build.latest.title"

$ws.Range("A283:C283").Copy()
$ws.Range("A284:C286").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B279").Select()
